$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to remain text so numeric-looking strings (e.g. "1.00", "0.532")
# are not auto-converted to numbers by Excel, matching the original inlineStr text cells.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '61.908.76'
$ws.Range("E2").Value = '  -0.68%  '

$ws.Range("D3").Value = '2.401.89'
$ws.Range("E3").Value = '  -1.03%  '

$ws.Range("E4").Value = '  -0.08%  '

$ws.Range("D5").Value = '561.95'
$ws.Range("E5").Value = '  +0.98%  '

$ws.Range("D6").Value = '142.29'
$ws.Range("E6").Value = '  -1.18%  '

$ws.Range("E7").Value = '  +0.00%  '

$ws.Range("D8").Value = '0.532'
$ws.Range("E8").Value = '  -0.94%  '

$ws.Range("D9").Value = '0.109'
$ws.Range("E9").Value = '  -0.61%  '

$ws.Range("E10").Value = '  -1.98%  '

$ws.Range("E11").Value = '  -2.60%  '

$ws.Range("E12").Value = '  -1.09%  '

$ws.Range("D13").Value = '25.52'
$ws.Range("E13").Value = '  -3.24%  '

$ws.Range("E14").Value = '  -1.83%  '

$ws.Range("D15").Value = '2.834.84'
$ws.Range("E15").Value = '  -1.21%  '

$ws.Range("D16").Value = '61.822.06'
$ws.Range("E16").Value = '  -0.59%  '

$ws.Range("D17").Value = '2.411.74'
$ws.Range("E17").Value = '  -0.73%  '

$ws.Range("E18").Value = '  +0.85%  '

$ws.Range("D19").Value = '321.42'
$ws.Range("E19").Value = '  -1.15%  '

$ws.Range("E20").Value = '  -1.10%  '

$ws.Range("E21").Value = '  +0.26%  '

$ws.Range("E22").Value = '  -0.38%  '

$ws.Range("D23").Value = '66.00'
$ws.Range("E23").Value = '  +1.40%  '

$ws.Range("E24").Value = '  -1.37%  '

$ws.Range("D25").Value = '8.79'
$ws.Range("E25").Value = '  -4.11%  '

$ws.Range("D26").Value = '560.02'
$ws.Range("E26").Value = '  -2.40%  '

$ws.Range("E27").Value = '  -0.60%  '

$ws.Range("D28").Value = '2.520.27'
$ws.Range("E28").Value = '  -0.71%  '

$ws.Range("D29").Value = '0.0₃0931'
$ws.Range("E29").Value = '  -1.51%  '

$ws.Range("D30").Value = '8.15'
$ws.Range("E30").Value = '  -2.90%  '

$ws.Range("E31").Value = '  -4.32%  '

$ws.Range("D32").Value = '0.147'
$ws.Range("E32").Value = '  -1.38%  '

$ws.Range("E33").Value = '  -0.09%  '

$ws.Range("E34").Value = '  -4.60%  '

$ws.Range("D35").Value = '1.00'
$ws.Range("E35").Value = '  +0.09%  '

$ws.Range("E36").Value = '  -3.70%  '

$ws.Range("B37").Value = 'Monero'
$ws.Range("C37").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D37").Value = '151.96'
$ws.Range("E37").Value = '  +2.83%  '

$ws.Range("B38").Value = 'RenderToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range("D38").Value = '5.41'
$ws.Range("E38").Value = '  -5.95%  '

$ws.Range("D39").Value = '0.377'
$ws.Range("E39").Value = '  -2.12%  '

$ws.Range("D40").Value = '18.56'
$ws.Range("E40").Value = '  -1.72%  '

$ws.Range("E41").Value = '  -6.86%  '

$ws.Range("E42").Value = '  -0.01%  '

$ws.Range("E43").Value = '  -3.45%  '

$ws.Range("D44").Value = '147.11'
$ws.Range("E44").Value = '  -3.14%  '

$ws.Range("E45").Value = '  -1.43%  '

$ws.Range("E46").Value = '  -3.26%  '

$ws.Range("D47").Value = '19.76'
$ws.Range("E47").Value = '  -3.59%  '

$ws.Range("E48").Value = '  -0.95%  '

$ws.Range("D49").Value = '0.0917'
$ws.Range("E49").Value = '  +0.27%  '

$ws.Range("E50").Value = '  -1.93%  '

$ws.Range("E51").Value = '  +0.48%  '
